# Applies the odds updates for Jogos_do_Dia_Betfair_Back_Lay_2026-01-05.xlsx
# as described by the source diff (commit: "Atualizando o arquivo XLSX").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("S2").Value = 3

# Row 3
$ws.Range("G3").Value = 2.24
$ws.Range("I3").Value = 5.1
$ws.Range("J3").Value = 3.25
$ws.Range("P3").Value = 1.98
$ws.Range("W3").Value = 1.8
$ws.Range("AK3").Value = 29

# Row 4
$ws.Range("P4").Value = 2.86

# Row 5
$ws.Range("F5").Value = 2.12
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 3.15
$ws.Range("L5").Value = 1.34
$ws.Range("M5").Value = 1.07
$ws.Range("P5").Value = 1.81
$ws.Range("Q5").Value = 2.02
$ws.Range("R5").Value = 1.28
$ws.Range("S5").Value = 3.6
$ws.Range("T5").Value = 1.78
$ws.Range("X5").Value = 15.5
$ws.Range("Y5").Value = 14
$ws.Range("Z5").Value = 29
$ws.Range("AA5").Value = 85
$ws.Range("AB5").Value = 9.6
$ws.Range("AC5").Value = 8.4
$ws.Range("AD5").Value = 17
$ws.Range("AE5").Value = 55
$ws.Range("AF5").Value = 15
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 20
$ws.Range("AI5").Value = 65
$ws.Range("AJ5").Value = 32
$ws.Range("AK5").Value = 27
$ws.Range("AN5").Value = 22

# Row 6
$ws.Range("F6").Value = 1.85
$ws.Range("H6").Value = 4.5
$ws.Range("I6").Value = 5.6
$ws.Range("K6").Value = 3.7
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 2.8
$ws.Range("P6").Value = 1.62
$ws.Range("Q6").Value = 2.3
$ws.Range("S6").Value = 4.6
$ws.Range("T6").Value = 2.08
$ws.Range("U6").Value = 1.76
$ws.Range("V6").Value = 1.22
$ws.Range("AC6").Value = 9.6

# Row 7
$ws.Range("F7").Value = 1.76
$ws.Range("G7").Value = 1.92
$ws.Range("J7").Value = 3.45
$ws.Range("K7").Value = 4.1
$ws.Range("M7").Value = 1.06
$ws.Range("Q7").Value = 1.91
$ws.Range("T7").Value = 1.81
$ws.Range("V7").Value = 1.18
$ws.Range("W7").Value = 2.08
$ws.Range("X7").Value = 18

# Row 8
$ws.Range("F8").Value = 2.28
$ws.Range("I8").Value = 3.65
$ws.Range("L8").Value = 1.46
$ws.Range("P8").Value = 1.74
$ws.Range("Q8").Value = 2.06
$ws.Range("U8").Value = 2
$ws.Range("AA8").Value = 70
$ws.Range("AJ8").Value = 36

# Row 9
$ws.Range("F9").Value = 1.89
$ws.Range("G9").Value = 2.14
$ws.Range("H9").Value = 4.2
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 3.2
$ws.Range("R9").Value = 1.28
$ws.Range("S9").Value = 3.75
$ws.Range("T9").Value = 1.87
$ws.Range("U9").Value = 1.9
$ws.Range("W9").Value = 1.88
$ws.Range("AB9").Value = 9.800000000000001
$ws.Range("AC9").Value = 9.800000000000001
$ws.Range("AG9").Value = 13
